# ------------------------------------------------------------------
# PlayerPerformance_4552.xlsx update
#   1. Remove the stray empty cell at B5 on "ODI Batting"
#   2. Add a new "ODI Batting Extra" sheet (after "ODI Bowling")
#      with MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 /
#      PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH columns
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Drop the empty inline-string cell left over at B5 ----------
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Range("B5").ClearContents()

# --- 2. Create "ODI Batting Extra" as the last sheet ----------------
$wsBowling = $wb.Worksheets.Item("ODI Bowling")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsBowling)
$newSheet.Name = "ODI Batting Extra"

# Re-use the bold/centered/bordered header style already used by the
# other sheets' header rows instead of creating a brand-new style.
$wsBatting.Range("A1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122) # xlPasteFormats

# --- Header row ------------------------------------------------------
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Data rows ---------------------------------------------------------
# MATCH_CODE, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL and MAN_OF_MATCH are
# stored as text (matching the rest of the workbook's export format);
# BATTING_POSITION is a genuine number. $null entries are left blank.
$rows = @(
    @{ Row = 2; MatchCode = "4454"; Position = 7;     Num4 = "7"; Num6 = "2"; Pct = "18.30%"; Mom = "NO" },
    @{ Row = 3; MatchCode = "4456"; Position = $null; Num4 = $null; Num6 = $null; Pct = $null; Mom = "NO" },
    @{ Row = 4; MatchCode = "4457"; Position = 7;     Num4 = "0"; Num6 = "0"; Pct = "7.60%";  Mom = "NO" },
    @{ Row = 5; MatchCode = "4480"; Position = 7;     Num4 = $null; Num6 = $null; Pct = $null; Mom = "NO" },
    @{ Row = 6; MatchCode = "4482"; Position = 7;     Num4 = "3"; Num6 = "0"; Pct = "12.64%"; Mom = "NO" }
)

foreach ($r in $rows) {
    $row = $r.Row

    $cellA = $newSheet.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.MatchCode
    $cellA.Style = "Normal"

    if ($null -ne $r.Position) {
        $newSheet.Cells.Item($row, 2).Value = $r.Position
    }

    if ($null -ne $r.Num4) {
        $cellC = $newSheet.Cells.Item($row, 3)
        $cellC.NumberFormat = "@"
        $cellC.Value = $r.Num4
        $cellC.Style = "Normal"
    }

    if ($null -ne $r.Num6) {
        $cellD = $newSheet.Cells.Item($row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $r.Num6
        $cellD.Style = "Normal"
    }

    if ($null -ne $r.Pct) {
        $cellE = $newSheet.Cells.Item($row, 5)
        $cellE.NumberFormat = "@"
        $cellE.Value = $r.Pct
        $cellE.Style = "Normal"
    }

    $cellF = $newSheet.Cells.Item($row, 6)
    $cellF.NumberFormat = "@"
    $cellF.Value = $r.Mom
    $cellF.Style = "Normal"
}
